$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# The "subgenus" field/column is being removed from the Materials sheet
# (and its corresponding "${subgenus}" mapping value below it) following
# the third round of review. Locate the header cell dynamically so the
# deletion is anchored on content rather than an assumed column letter.
$headerCell = $ws.Rows.Item(1).Find("subgenus")
if ($headerCell -ne $null) {
    $headerCell.EntireColumn.Delete()
}
